# HOLON Base/db_backboneConfig.xlsx update
# - District heating network: heat nodes, DistrictHeating netConnection agent,
#   House/Building netConnection agents may reference both an electrical
#   parent node and a heat parent node.
# - Minor index fix on the default storage energy assets sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) config_netNodes  -- cosmetic only (column D width + last selection)
# ---------------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("config_netNodes")
$wsNodes.Columns.Item(4).ColumnWidth = 12.428571428571429

# ---------------------------------------------------------------------------
# 2) config_netConnections -- the main data change
#    Old layout: A index | B agenttype | C id | D type | E type2 | F parent | G capacity_kw
#    New layout: A index | B agenttype | C id | D type | E type2 |
#                F parent_electric | G parent_heat | H capacity_kw
# ---------------------------------------------------------------------------
$wsConn = $wb.Worksheets.Item("config_netConnections")

# -- New row 29: District Heating netConnection agent. The new shared
#    strings are introduced in this exact order (DISTRICTHEATING, b28)
#    to match how the workbook's sharedStrings table was actually grown --
$wsConn.Cells.Item(29, 4).Value2 = "DISTRICTHEATING"
$wsConn.Cells.Item(29, 3).Value2 = "b28"
$wsConn.Cells.Item(29, 1).Value2 = 27
$wsConn.Cells.Item(29, 2).Value2 = "netConnection"
$wsConn.Cells.Item(29, 6).Value2 = "E2"
$wsConn.Cells.Item(29, 7).Value2 = "H1"
$wsConn.Cells.Item(29, 8).Value2 = 400000

# -- Header row: rename F1, insert new G1, move old G1 (capacity_kw) to H1 --
$wsConn.Cells.Item(1, 6).Value2 = "parent_electric"
$wsConn.Cells.Item(1, 7).Clear()
$wsConn.Cells.Item(1, 7).Value2 = "parent_heat"
$wsConn.Cells.Item(1, 8).Value2 = "capacity_kw"
$wsConn.Cells.Item(1, 8).NumberFormat = "0"

# -- Shift the capacity_kw values (old col G) from G to H for every data row,
#    then clear the vacated G column so it can be reused for parent_heat --
for ($r = 2; $r -le 28; $r++) {
    $srcCell = $wsConn.Cells.Item($r, 7)
    $dstCell = $wsConn.Cells.Item($r, 8)
    if ($srcCell.HasFormula) {
        $dstCell.Formula = $srcCell.Formula
    } else {
        $dstCell.Value2 = $srcCell.Value2
    }
    $dstCell.NumberFormat = "0"
}
$wsConn.Range("G2:G28").Clear()

# -- House / Building netConnection agents that also connect to the heat
#    node (H1) get a parent_heat entry in column G --
$heatConnectedRows = 4, 11, 12, 13, 20
foreach ($r in $heatConnectedRows) {
    $wsConn.Cells.Item($r, 7).Value2 = "H1"
}

# -- Column widths (best effort; engine quantizes to whole pixels) --
$wsConn.Columns.Item(4).ColumnWidth = 17.285714285714285
$wsConn.Columns.Item(5).ColumnWidth = 16.571428571428573
$wsConn.Columns.Item(6).ColumnWidth = 13.428571428571429
$wsConn.Columns.Item(7).ColumnWidth = 11.857142857142858

# ---------------------------------------------------------------------------
# 3) config_energyAssets -- fix index numbering on the last two storage rows
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("config_energyAssets")
$wsAssets.Cells.Item(6, 1).Value2 = 4
$wsAssets.Cells.Item(7, 1).Value2 = 5

# ---------------------------------------------------------------------------
# 4) Selections -- replicate the last-saved cursor position per sheet.
#    Doing this last-to-first in tab order makes config_energyAssets end up
#    as the active sheet again (it was, and stays, the active tab).
# ---------------------------------------------------------------------------
$wsNodes.Range("E5").Select() | Out-Null
$wsConn.Range("D29").Select() | Out-Null
$wsAssets.Range("J21").Select() | Out-Null
